# localization-status.xlsx: "Generate Report for Archive"
#
# 1) The localization status of the two handed-off files moves from
#    "Ready for handoff" to "In Translation" on every sheet that shows it
#    (the Overview sheet's per-language status columns, plus the per-locale
#    "Status" column on the zh-cn and de-de sheets).
# 2) Because the new status text is shorter, the status column(s) are
#    re-sized (narrower) to fit the content, on all three sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# The width the "Status" columns are fit to (their content/format no longer
# needs the wider "Ready for handoff" string).
$newStatusColumnWidth = 12.576851254417766

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = $newStatus
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $newStatusColumnWidth

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = $newStatus
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newStatusColumnWidth

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = $newStatus
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newStatusColumnWidth
